# This workbook contains NATMI ligand-receptor pair output recomputed with
# new TPM values. The underlying per-cluster ligand/receptor average and
# total expression values changed for the "ECs" cluster, which cascades
# into every derived specificity / edge-weight column for all 9 data rows.
# We set every affected cell to its freshly recomputed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=ECs, Target=ECs)
$ws.Range("G2").Value = 26.92947733333333
$ws.Range("H2").Value = 80.788432
$ws.Range("I2").Value = 0.9279949792877585
$ws.Range("J2").Value = 0.9279949792877585
$ws.Range("M2").Value = 8.119312000000001
$ws.Range("N2").Value = 24.357936
$ws.Range("O2").Value = 0.6660224804915575
$ws.Range("P2").Value = 0.6660224804915574
$ws.Range("Q2").Value = 218.6488284662614
$ws.Range("R2").Value = 1967.839456196352
$ws.Range("S2").Value = 0.6180655179889445
$ws.Range("T2").Value = 0.6180655179889444

# Row 3 (Sending=ECs, Target=FAPs)
$ws.Range("G3").Value = 26.92947733333333
$ws.Range("H3").Value = 80.788432
$ws.Range("I3").Value = 0.9279949792877585
$ws.Range("J3").Value = 0.9279949792877585
$ws.Range("O3").Value = 0.2728595563580828
$ws.Range("P3").Value = 0.2728595563580827
$ws.Range("Q3").Value = 89.57719008146132
$ws.Range("R3").Value = 806.1947107331519
$ws.Range("S3").Value = 0.253212298350986
$ws.Range("T3").Value = 0.253212298350986

# Row 4 (Sending=ECs, Target=MuSCs)
$ws.Range("G4").Value = 26.92947733333333
$ws.Range("H4").Value = 80.788432
$ws.Range("I4").Value = 0.9279949792877585
$ws.Range("J4").Value = 0.9279949792877585
$ws.Range("O4").Value = 0.06111796315035969
$ws.Range("P4").Value = 0.06111796315035968
$ws.Range("Q4").Value = 20.06444441816356
$ws.Range("R4").Value = 180.579999763472
$ws.Range("S4").Value = 0.05671716294782803
$ws.Range("T4").Value = 0.05671716294782802

# Row 5 (Sending=FAPs, Target=ECs)
$ws.Range("I5").Value = 0.04103565698374688
$ws.Range("J5").Value = 0.04103565698374688
$ws.Range("M5").Value = 8.119312000000001
$ws.Range("N5").Value = 24.357936
$ws.Range("O5").Value = 0.6660224804915575
$ws.Range("P5").Value = 0.6660224804915574
$ws.Range("Q5").Value = 9.668584987093334
$ws.Range("R5").Value = 87.01726488384
$ws.Range("S5").Value = 0.0273306700529158
$ws.Range("T5").Value = 0.0273306700529158

# Row 6 (Sending=FAPs, Target=FAPs)
$ws.Range("I6").Value = 0.04103565698374688
$ws.Range("J6").Value = 0.04103565698374688
$ws.Range("O6").Value = 0.2728595563580828
$ws.Range("P6").Value = 0.2728595563580827
$ws.Range("S6").Value = 0.01119697115944764
$ws.Range("T6").Value = 0.01119697115944763

# Row 7 (Sending=FAPs, Target=MuSCs)
$ws.Range("I7").Value = 0.04103565698374688
$ws.Range("J7").Value = 0.04103565698374688
$ws.Range("O7").Value = 0.06111796315035969
$ws.Range("P7").Value = 0.06111796315035968
$ws.Range("S7").Value = 0.002508015771383442
$ws.Range("T7").Value = 0.002508015771383442

# Row 8 (Sending=MuSCs, Target=ECs)
$ws.Range("G8").Value = 0.8986996666666666
$ws.Range("I8").Value = 0.03096936372849452
$ws.Range("J8").Value = 0.03096936372849452
$ws.Range("M8").Value = 8.119312000000001
$ws.Range("N8").Value = 24.357936
$ws.Range("O8").Value = 0.6660224804915575
$ws.Range("P8").Value = 0.6660224804915574
$ws.Range("Q8").Value = 7.296822987962667
$ws.Range("R8").Value = 65.67140689166401
$ws.Range("S8").Value = 0.02062629244969719
$ws.Range("T8").Value = 0.02062629244969719

# Row 9 (Sending=MuSCs, Target=FAPs)
$ws.Range("G9").Value = 0.8986996666666666
$ws.Range("I9").Value = 0.03096936372849452
$ws.Range("J9").Value = 0.03096936372849452
$ws.Range("O9").Value = 0.2728595563580828
$ws.Range("P9").Value = 0.2728595563580827
$ws.Range("Q9").Value = 2.989400420612666
$ws.Range("S9").Value = 0.008450286847649115
$ws.Range("T9").Value = 0.008450286847649115

# Row 10 (Sending=MuSCs, Target=MuSCs)
$ws.Range("G10").Value = 0.8986996666666666
$ws.Range("I10").Value = 0.03096936372849452
$ws.Range("J10").Value = 0.03096936372849452
$ws.Range("O10").Value = 0.06111796315035969
$ws.Range("P10").Value = 0.06111796315035968
$ws.Range("Q10").Value = 0.6695974558754444
$ws.Range("R10").Value = 6.026377102879
$ws.Range("S10").Value = 0.001892784431148214
$ws.Range("T10").Value = 0.001892784431148214
